$wb = $excel.ActiveWorkbook

# Delete the empty "Sheet2" worksheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()

# Rename "Sheet1" (the worksheet holding the task data) to the new name
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "taskflow -frontEnd"

# Set a shared string value of 5 spaces into E28
$ws.Range("E28").Value = "     "

# Widen column A (target stored width ~11.7265625 chars; closest value the
# engine's column-width quantization can reproduce)
$ws.Columns.Item(1).ColumnWidth = 10.8125

# Update the selection on the sheet
$ws.Range("F4").Select()
